$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain text so values are not
# auto-converted to numbers by Excel (matches original inline-string cells)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.315.12'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '1.879.38'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '237.37'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = '0.4818'
$ws.Range('E7').Value = '  -2.31%  '
$ws.Range('D8').Value = '0.2888'
$ws.Range('E8').Value = '  -3.12%  '
$ws.Range('D9').Value = '0.06584'
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').Value = '1.882.29'
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('D11').Value = '16.93'
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').Value = '0.07388'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '5.190'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').Value = '87.98'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '0.6599'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '30.279.87'
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('D17').Value = '13.61'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '0.9998'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '0.000007725'
$ws.Range('E19').Value = '  -2.78%  '
$ws.Range('D20').Value = '5.466'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').Value = '2.138.39'
$ws.Range('E21').Value = '  -0.52%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = '194.82'
$ws.Range('E23').Value = '  -4.01%  '
$ws.Range('D24').Value = '6.168'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('D26').Value = '164.33'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = '18.25'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('D28').Value = '1.927'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').Value = '1.441'
$ws.Range('D30').Value = '4.276'
$ws.Range('D31').Value = '0.09140'
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').Value = '4.052'
$ws.Range('D33').Value = '0.05055'
$ws.Range('E33').Value = '  -4.21%  '
$ws.Range('D34').Value = '0.7425'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +1.83%  '
$ws.Range('D36').Value = '2.710'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').Value = '0.01842'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').Value = '2.075'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '106.14'
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.4323'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.874'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('D44').Value = '0.9992'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = '7.640'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('D47').Value = '1.582'
$ws.Range('E47').Value = '  +10.03%  '
$ws.Range('D48').Value = '65.17'
$ws.Range('E48').Value = '  -11.18%  '
$ws.Range('D49').Value = '8.892'
$ws.Range('E49').Value = '  -1.40%  '
$ws.Range('D50').Value = '34.17'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('D51').Value = '0.05721'
$ws.Range('E51').Value = '  -2.74%  '
